$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "45.547.34"
Set-TextValue "E2" "  -2.07%  "
Set-TextValue "D3" "2.400.20"
Set-TextValue "E3" "  +4.52%  "
Set-TextValue "E4" "  +0.11%  "
Set-TextValue "D5" "299.81"
Set-TextValue "E5" "  -1.57%  "
Set-TextValue "D6" "97.49"
Set-TextValue "E6" "  -3.69%  "
Set-TextValue "D7" "0.563"
Set-TextValue "E7" "  -0.76%  "
Set-TextValue "E8" "  +0.10%  "
Set-TextValue "D9" "0.511"
Set-TextValue "E9" "  -2.43%  "
Set-TextValue "E10" "  -4.16%  "
Set-TextValue "D11" "0.0792"
Set-TextValue "E11" "  +0.51%  "
Set-TextValue "D12" "7.16"
Set-TextValue "E12" "  -3.44%  "
Set-TextValue "E13" "  +1.00%  "
Set-TextValue "D14" "2.754.02"
Set-TextValue "E14" "  +4.21%  "
Set-TextValue "D15" "2.408.70"
Set-TextValue "E15" "  +5.14%  "
Set-TextValue "D16" "0.844"
Set-TextValue "E16" "  +3.74%  "
Set-TextValue "D17" "14.23"
Set-TextValue "E17" "  +2.82%  "
Set-TextValue "D18" "45.500.23"
Set-TextValue "E18" "  -2.10%  "
Set-TextValue "D19" "12.95"
Set-TextValue "E19" "  -1.17%  "
Set-TextValue "D20" "0.0₃0950"
Set-TextValue "E20" "  +1.21%  "
Set-TextValue "D21" "6.23"
Set-TextValue "E21" "  +3.43%  "
Set-TextValue "D22" "67.12"
Set-TextValue "E22" "  +1.47%  "
Set-TextValue "D23" "242.88"
Set-TextValue "E23" "  -2.44%  "
Set-TextValue "D24" "2.82"
Set-TextValue "E24" "  -2.44%  "
Set-TextValue "E25" "  -0.03%  "
Set-TextValue "D26" "1.93"
Set-TextValue "E26" "  +0.13%  "
Set-TextValue "D27" "38.38"
Set-TextValue "E27" "  -9.99%  "
Set-TextValue "E28" "  -1.96%  "
Set-TextValue "D29" "9.78"
Set-TextValue "E29" "  -1.13%  "
Set-TextValue "D30" "3.84"
Set-TextValue "E30" "  +17.26%  "
Set-TextValue "D31" "21.20"
Set-TextValue "E31" "  +5.86%  "
Set-TextValue "D32" "2.74"
Set-TextValue "E32" "  -3.01%  "
Set-TextValue "B33" "Monero"
Set-TextValue "C33" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D33" "148.86"
Set-TextValue "E33" "  +0.76%  "
Set-TextValue "B34" "Filecoin"
Set-TextValue "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "5.52"
Set-TextValue "E34" "  -2.25%  "
Set-TextValue "D35" "0.0775"
Set-TextValue "E35" "  -2.52%  "
Set-TextValue "B36" "Kaspa"
Set-TextValue "C36" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.114"
Set-TextValue "E36" "  +0.23%  "
Set-TextValue "B37" "ARBITRUM"
Set-TextValue "C37" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D37" "1.97"
Set-TextValue "E37" "  +11.33%  "
Set-TextValue "E38" "  -1.28%  "
Set-TextValue "D39" "15.27"
Set-TextValue "E39" "  -4.81%  "
Set-TextValue "D40" "3.86"
Set-TextValue "E40" "  -3.57%  "
Set-TextValue "D41" "0.0300"
Set-TextValue "E41" "  -0.65%  "
Set-TextValue "D42" "3.27"
Set-TextValue "E42" "  -2.42%  "
Set-TextValue "D43" "1.944.15"
Set-TextValue "E43" "  +6.72%  "
Set-TextValue "D44" "1.00"
Set-TextValue "E44" "  +0.12%  "
Set-TextValue "D45" "91.89"
Set-TextValue "E45" "  +4.12%  "
Set-TextValue "D46" "1.75"
Set-TextValue "E46" "  -10.74%  "
Set-TextValue "D47" "8.72"
Set-TextValue "E47" "  +10.40%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "15.55"
Set-TextValue "E48" "  +16.33%  "
Set-TextValue "B49" "Aave"
Set-TextValue "C49" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D49" "103.09"
Set-TextValue "E49" "  +7.40%  "
Set-TextValue "E50" "  -3.26%  "
Set-TextValue "D51" "2.638.85"
Set-TextValue "E51" "  +4.73%  "
